$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.873.24'
$ws.Range("E2").Value = '  +7.74%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.815.90'
$ws.Range("E3").Value = '  +5.20%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9983'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.42'
$ws.Range("E5").Value = '  +4.22%  '

$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9982'
$ws.Range("E6").Value = '  -0.02%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4988'
$ws.Range("E7").Value = '  +2.52%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2798'
$ws.Range("E8").Value = '  +8.26%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06401'
$ws.Range("E9").Value = '  +3.23%  '

$ws.Range("B10").Value = 'WrappedEther'
$ws.Range("C10").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.815.89'
$ws.Range("E10").Value = '  +5.19%  '

$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.86'
$ws.Range("E11").Value = '  +5.80%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07168'
$ws.Range("E12").Value = '  +3.81%  '

$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6520'
$ws.Range("E13").Value = '  +7.64%  '

$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '83.14'
$ws.Range("E14").Value = '  +8.11%  '

$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.729'
$ws.Range("E15").Value = '  +5.59%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.834.98'
$ws.Range("E16").Value = '  +8.41%  '

$ws.Range("B17").Value = 'Dai'
$ws.Range("C17").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9994'
$ws.Range("E17").Value = '  +0.10%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007403'
$ws.Range("E18").Value = '  +3.32%  '

$ws.Range("B19").Value = 'BinanceUSD'
$ws.Range("C19").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9978'
$ws.Range("E19").Value = '  -0.04%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.32'
$ws.Range("E20").Value = '  +7.77%  '

$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.050.25'
$ws.Range("E21").Value = '  +5.04%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.627'
$ws.Range("E22").Value = '  +4.78%  '

$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.914'
$ws.Range("E23").Value = '  +4.15%  '

$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.364'
$ws.Range("E24").Value = '  +5.79%  '

$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.14'
$ws.Range("E25").Value = '  +4.85%  '

$ws.Range("B26").Value = 'BitcoinCash'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '127.88'
$ws.Range("E26").Value = '  +20.58%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.30'
$ws.Range("E27").Value = '  +7.01%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.905'
$ws.Range("E28").Value = '  +7.32%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.396'
$ws.Range("E29").Value = '  +1.22%  '

$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.195'
$ws.Range("E30").Value = '  +6.76%  '

$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08389'
$ws.Range("E31").Value = '  +5.30%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.874'
$ws.Range("E32").Value = '  +5.39%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04987'
$ws.Range("E33").Value = '  +10.73%  '

$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.095'
$ws.Range("E34").Value = '  +8.47%  '

$ws.Range("B35").Value = 'Frax'
$ws.Range("C35").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9967'
$ws.Range("E35").Value = '  -0.09%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6830'
$ws.Range("E36").Value = '  +9.53%  '

$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.683'
$ws.Range("E37").Value = '  +3.18%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.744'
$ws.Range("E38").Value = '  +13.13%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.221'
$ws.Range("E39").Value = '  +9.39%  '

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9710'
$ws.Range("E40").Value = '  +4.66%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01595'
$ws.Range("E41").Value = '  +6.57%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.018'
$ws.Range("E42").Value = '  +6.38%  '

$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9989'
$ws.Range("E43").Value = '  +0.12%  '

$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.74'
$ws.Range("E44").Value = '  +2.49%  '

$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4133'
$ws.Range("E45").Value = '  +7.91%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.266'
$ws.Range("E46").Value = '  +6.09%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1230'
$ws.Range("E47").Value = '  +6.25%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05512'
$ws.Range("E48").Value = '  +2.14%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.211'
$ws.Range("E49").Value = '  +3.83%  '

$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '31.75'
$ws.Range("E50").Value = '  +5.49%  '

$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3660'
$ws.Range("E51").Value = '  +8.68%  '

